# Update page and figure
#
# Shrinks the font used in the three "CaixaDeTexto" label boxes
# (cores_de_pele, alturas, and the "Tabela resultado do JOIN" caption) from
# 20pt Inconsolata to 16pt Nunito Sans Black, and repositions/resizes the
# boxes to match the new, smaller auto-fit text (the third box is also
# moved and narrowed, which makes its text wrap onto more lines so it grows
# taller even though the font got smaller).

# Shape.Left/Top/Width/Height take points but are stored as EMU (1 pt =
# 12700 EMU) after being narrowed to a 32-bit float; EmuToPtPrecise finds a
# point value whose float32 truncation, scaled back up, lands exactly back
# on the desired EMU integer so the saved XML matches to the EMU.
function EmuToPtPrecise([double]$targetEmu) {
    $pt = $targetEmu / 12700.0
    $step = 0.0000001
    for ($i = 0; $i -lt 200000; $i++) {
        $f = [float]$pt
        $got = [math]::floor($f * 12700)
        if ($got -eq $targetEmu) {
            return $pt
        }
        $pt = $pt + $step
    }
    return $targetEmu / 12700.0
}

function Set-LabelFont($shape) {
    $tr = $shape.TextFrame.TextRange
    $tr.Font.Size = 16
    $tr.Font.Name = "Nunito Sans Black"
}

function Set-ShapeRectEmu($shape, $offX, $offY, $cx, $cy) {
    $shape.Left   = EmuToPtPrecise $offX
    $shape.Top    = EmuToPtPrecise $offY
    $shape.Width  = EmuToPtPrecise $cx
    $shape.Height = EmuToPtPrecise $cy
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "cores_de_pele" textbox -------------------------------------------
$sh1 = $s.Shapes.Item("CaixaDeTexto 15")
Set-LabelFont $sh1
Set-ShapeRectEmu $sh1 648476 525978 1924665 338554

# --- "alturas" textbox --------------------------------------------------
$sh2 = $s.Shapes.Item("CaixaDeTexto 17")
Set-LabelFont $sh2
Set-ShapeRectEmu $sh2 4077978 522889 1924665 338554

# --- "Tabela resultado do JOIN" textbox ---------------------------------
$sh3 = $s.Shapes.Item("CaixaDeTexto 20")
Set-LabelFont $sh3
Set-ShapeRectEmu $sh3 7172238 358247 2023945 584775
